# Daily attendance processing - 2025-11-26 20:27:36
# Normalises the "Recorded By" (column G) cell text so that, for rows
# recording exactly two recorders, the human account is listed before
# "System" (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
# Rows whose recorder list includes the backup/service account
# (backup@backdoor.com) or that don't have exactly two comma-separated
# names are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value()

    if ($null -eq $value) { continue }
    if ($value -notlike "*,*") { continue }
    if ($value -like "*backup@backdoor.com*") { continue }

    $parts = $value -split ",\s*"
    if ($parts.Count -ne 2) { continue }

    $swapped = "{0}, {1}" -f $parts[1].Trim(), $parts[0].Trim()
    if ($swapped -ne $value) {
        $cell.Value = $swapped
    }
}
